$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 1341.7826
$ws.Range("I28").Value = 528.2
$ws.Range("J28").Value = 2867.25
$ws.Range("K28").Value = 528.2
$ws.Range("L28").Value = 2867.25
$ws.Range("M28").Value = -43.20000000000005
$ws.Range("N28").Value = -3837.25
# Row 40
$ws.Range("H40").Value = 4324.3125
$ws.Range("I40").Value = 5957
$ws.Range("J40").Value = 3582.182
$ws.Range("K40").Value = 5957
$ws.Range("L40").Value = 3582.182
$ws.Range("M40").Value = -5782
$ws.Range("N40").Value = -3932.182
# Row 98
$ws.Range("H98").Value = 1961.878
$ws.Range("I98").Value = 1965.925
$ws.Range("K98").Value = 1965.925
$ws.Range("M98").Value = -467.925
# Row 103
$ws.Range("H103").Value = 3095.1428
$ws.Range("J103").Value = 1694.5
$ws.Range("L103").Value = 5083.5
$ws.Range("N103").Value = -6255.5
# Row 107
$ws.Range("H107").Value = 1098.4546
$ws.Range("I107").Value = 1309.7142
$ws.Range("J107").Value = 728.75
$ws.Range("K107").Value = 1309.7142
$ws.Range("L107").Value = 728.75
$ws.Range("M107").Value = 610.2858000000001
$ws.Range("N107").Value = -4568.75
# Row 122
$ws.Range("H122").Value = 1961.878
$ws.Range("I122").Value = 1965.925
$ws.Range("K122").Value = 5897.775
$ws.Range("M122").Value = -3447.775
# Row 134
$ws.Range("H134").Value = 68761.42999999999
$ws.Range("J134").Value = 68761.42999999999
$ws.Range("L134").Value = 68761.42999999999
$ws.Range("N134").Value = -78901.42999999999
# Row 141
$ws.Range("H141").Value = 5479.125
$ws.Range("I141").Value = 3250
$ws.Range("K141").Value = 9750
$ws.Range("M141").Value = -4570

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 207606.14
$ws.Range("I74").Value = 278832.7
$ws.Range("J74").Value = 4101.7144
$ws.Range("K74").Value = 278832.7
$ws.Range("L74").Value = 4101.7144
$ws.Range("M74").Value = -277958.7
$ws.Range("N74").Value = -5849.7144
# Row 77
$ws.Range("H77").Value = 207606.14
$ws.Range("I77").Value = 278832.7
$ws.Range("J77").Value = 4101.7144
$ws.Range("K77").Value = 1394163.5
$ws.Range("L77").Value = 20508.572
$ws.Range("M77").Value = -1389795.5
$ws.Range("N77").Value = -29244.572
# Row 88
$ws.Range("H88").Value = 3871.9092
$ws.Range("I88").Value = 1772
$ws.Range("K88").Value = 1772
$ws.Range("M88").Value = -1366
# Row 91
$ws.Range("H91").Value = 3871.9092
$ws.Range("I91").Value = 1772
$ws.Range("K91").Value = 1772
$ws.Range("M91").Value = -368
# Row 122
$ws.Range("H122").Value = 5622.0625
$ws.Range("I122").Value = 4992.4
$ws.Range("K122").Value = 14977.2
$ws.Range("M122").Value = -12527.2

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3483.2222
$ws.Range("I86").Value = 3106.125
$ws.Range("K86").Value = 3106.125
$ws.Range("M86").Value = -1983.125
# Row 89
$ws.Range("H89").Value = 3483.2222
$ws.Range("I89").Value = 3106.125
$ws.Range("K89").Value = 15530.625
$ws.Range("M89").Value = -9914.625
# Row 139
$ws.Range("H139").Value = 81272.42999999999
$ws.Range("J139").Value = 81272.42999999999
$ws.Range("L139").Value = 81272.42999999999
$ws.Range("N139").Value = -91552.42999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Range("H52").Value = 95809
$ws.Range("J52").Value = 95809
$ws.Range("L52").Value = 95809
$ws.Range("N52").Value = -96397
# Row 55
$ws.Range("H55").Value = 14124.75
$ws.Range("J55").Value = 14124.75
$ws.Range("L55").Value = 14124.75
$ws.Range("N55").Value = -14754.75
# Row 116
$ws.Range("H116").Value = 64669.25
$ws.Range("J116").Value = 64669.25
$ws.Range("L116").Value = 64669.25
$ws.Range("N116").Value = -73847.25
# Row 132
$ws.Range("H132").Value = 3114.2812
$ws.Range("I132").Value = 3143.5454
$ws.Range("J132").Value = 3049.9
$ws.Range("K132").Value = 9430.636200000001
$ws.Range("L132").Value = 9149.700000000001
$ws.Range("M132").Value = -6900.636200000001
$ws.Range("N132").Value = -14209.7
# Row 134
$ws.Range("H134").Value = 3069.3142
$ws.Range("I134").Value = 2982.9614
$ws.Range("J134").Value = 3318.7778
$ws.Range("K134").Value = 8948.8842
$ws.Range("L134").Value = 9956.3334
$ws.Range("M134").Value = -6413.8842
$ws.Range("N134").Value = -15026.3334
# Row 138
$ws.Range("H138").Value = 91922.30499999999
$ws.Range("J138").Value = 91922.30499999999
$ws.Range("L138").Value = 91922.30499999999
$ws.Range("N138").Value = -102202.305

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 800
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 53
$ws.Range("H53").Value = 800
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
# Row 81
$ws.Range("H81").Value = 2479.8
$ws.Range("I81").Value = 800
$ws.Range("J81").Value = 4999.5
$ws.Range("K81").Value = 2400
$ws.Range("L81").Value = 14998.5
$ws.Range("M81").Value = -1277
$ws.Range("N81").Value = -17244.5
# Row 84
$ws.Range("H84").Value = 2479.8
$ws.Range("I84").Value = 800
$ws.Range("J84").Value = 4999.5
$ws.Range("K84").Value = 7200
$ws.Range("L84").Value = 44995.5
$ws.Range("M84").Value = -1584
$ws.Range("N84").Value = -56227.5
# Row 131
$ws.Range("H131").Value = 14548.454
$ws.Range("J131").Value = 2134
$ws.Range("L131").Value = 6402
$ws.Range("N131").Value = -16482
# Row 132
$ws.Range("H132").Value = 2593.05
$ws.Range("I132").Value = 1485.7
$ws.Range("K132").Value = 13371.3
$ws.Range("M132").Value = -10841.3

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 8231.286
$ws.Range("I102").Value = 572.4
$ws.Range("K102").Value = 572.4
$ws.Range("M102").Value = 1049.6
# Row 113
$ws.Range("H113").Value = 3798.2
$ws.Range("I113").Value = 3107.6667
$ws.Range("K113").Value = 3107.6667
$ws.Range("M113").Value = -937.6667000000002

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 366307
$ws.Range("J94").Value = 366307
$ws.Range("L94").Value = 366307
$ws.Range("N94").Value = -367659
# Row 132
$ws.Range("H132").Value = 6391.923
$ws.Range("I132").Value = 3011.5881
$ws.Range("K132").Value = 9034.764299999999
$ws.Range("M132").Value = -6504.764299999999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 4799.5386
$ws.Range("I96").Value = 5049.625
$ws.Range("J96").Value = 4399.4
$ws.Range("K96").Value = 5049.625
$ws.Range("L96").Value = 4399.4
$ws.Range("M96").Value = -3676.625
$ws.Range("N96").Value = -7145.4
# Row 122
$ws.Range("H122").Value = 7814504.5
$ws.Range("I122").Value = 1952.6086
$ws.Range("K122").Value = 5857.825800000001
$ws.Range("M122").Value = -3407.825800000001
# Row 138
$ws.Range("H138").Value = 85398.8
$ws.Range("J138").Value = 85398.8
$ws.Range("L138").Value = 85398.8
$ws.Range("N138").Value = -95678.8
